# Decrement the "剩余" (remaining) value in column E for each data row
# (rows 2-99) by 1, except for row 36 whose F (start date) value is not a
# valid date and was therefore skipped by the original automated update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$skipRows = @(36)

for ($r = 2; $r -le 99; $r++) {
    if ($skipRows -contains $r) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)  # column E = 5
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current - 1
    }
}
